$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 497/498; existing rows 497-538 shift down to 499-540
$ws.Rows("497:498").Insert()

# New row 497 - Primera, weekly update for 2022-01-17 (serial 44578)
$ws.Range("A497").Value = 9
$ws.Range("B497").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C497").Value = "Metropolitana"
$ws.Range("D497").Value = 44578
$ws.Range("E497").Value = 13
$ws.Range("F497").Value = 100112023
$ws.Range("G497").Value = "Brócoli"
$ws.Range("H497").Value = "Sin especificar"
$ws.Range("I497").Value = "Primera"
$ws.Range("J497").Value = 1600
$ws.Range("K497").Value = 850
$ws.Range("L497").Value = 900
$ws.Range("M497").Value = 875
$ws.Range("N497").Value = "`$/unidad"
$ws.Range("O497").Value = "Región Metropolitana"
$ws.Range("P497").Value = 875
$ws.Range("Q497").Value = 1
$ws.Range("R497").Value = "Hortaliza"

# New row 498 - Segunda, same date
$ws.Range("A498").Value = 9
$ws.Range("B498").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C498").Value = "Metropolitana"
$ws.Range("D498").Value = 44578
$ws.Range("E498").Value = 13
$ws.Range("F498").Value = 100112023
$ws.Range("G498").Value = "Brócoli"
$ws.Range("H498").Value = "Sin especificar"
$ws.Range("I498").Value = "Segunda"
$ws.Range("J498").Value = 790
$ws.Range("K498").Value = 700
$ws.Range("L498").Value = 750
$ws.Range("M498").Value = 725
$ws.Range("N498").Value = "`$/unidad"
$ws.Range("O498").Value = "Región Metropolitana"
$ws.Range("P498").Value = 725
$ws.Range("Q498").Value = 1
$ws.Range("R498").Value = "Hortaliza"
